# Apply the four content fixes described in the commit:
# "BienRegistrar Fix se graba el costo ahora."
#
# 1) Correct the date shown at the top of the letter.
# 2) Correct the Fiscalia number referenced in the body.
# 3) Correct the requested amount ("el costo").
# 4) Fill in the justification text that was left blank/with a single space.

$d = $word.ActiveDocument

# 1. Fecha: 01 -> 09 de diciembre de 2017.
$d.Content.Find.Execute("01 de diciembre de 2017.", $false, $false, $false, $false, $false, `
    $true, 1, $false, "09 de diciembre de 2017.", 2)

# 2. Fiscalia Nacional en lo Criminal y Correccional Nro 3 -> Nro 2
$d.Content.Find.Execute("Fiscalia Nacional en lo Criminal y Correccional Nro 3", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Fiscalia Nacional en lo Criminal y Correccional Nro 2", 2)

# 3. Monto solicitado: $ 500,00 -> $ 1.000,00
$d.Content.Find.Execute("$ 500,00", $false, $false, $false, $false, $false, `
    $true, 1, $false, "$ 1.000,00", 2)

# 4. Justificacion field result: " " -> justification text.
# Target the specific DOCPROPERTY field (PJustificacion) instead of a blind
# text search, since several other runs in the document also contain a
# single space.
foreach ($f in $d.Fields) {
    if ($f.Code.Text -match "PJustificacion") {
        $result = $f.Result
        $rng = $d.Range($result.Start, $result.End)
        $rng.Text = "Finalmente, la presente erogación de fondos es solicitada por este curso debido a que Justif"
    }
}
